$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text value on a cell while keeping it stored as Text and
# preventing Excel auto-converting numeric-looking strings (e.g. "244.40")
# into a Number (which would drop significant trailing zeros). We briefly
# force the Text number format, assign the value, then restore the cell to
# the "Normal" style so no stray formatting is left behind.
function Set-TextValue($rangeAddr, $val) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '29.230.43'
$ws.Range("E2").Value = '  +0.47%  '

$ws.Range("D3").Value = '1.844.12'
$ws.Range("E3").Value = '  +0.41%  '

$ws.Range("E4").Value = '  +0.06%  '

Set-TextValue "D5" '244.40'
$ws.Range("E5").Value = '  -0.16%  '

Set-TextValue "D6" '0.6282'
$ws.Range("E6").Value = '  -1.48%  '

$ws.Range("E7").Value = '  +0.15%  '

Set-TextValue "D8" '0.07538'
$ws.Range("E8").Value = '  -0.40%  '

Set-TextValue "D9" '0.2951'
$ws.Range("E9").Value = '  -0.10%  '

Set-TextValue "D10" '23.38'
$ws.Range("E10").Value = '  +1.88%  '

Set-TextValue "D11" '0.07726'
$ws.Range("E11").Value = '  -0.41%  '

$ws.Range("D12").Value = '1.850.47'
$ws.Range("E12").Value = '  +0.26%  '

Set-TextValue "D13" '5.029'
$ws.Range("E13").Value = '  +0.22%  '

Set-TextValue "D14" '0.6790'
$ws.Range("E14").Value = '  +0.85%  '

Set-TextValue "D15" '83.33'
$ws.Range("E15").Value = '  -0.05%  '

Set-TextValue "D16" '0.000009274'
$ws.Range("E16").Value = '  -2.83%  '

Set-TextValue "D17" '5.990'
$ws.Range("E17").Value = '  -2.24%  '

$ws.Range("D18").Value = '29.218.84'
$ws.Range("E18").Value = '  +0.35%  '

$ws.Range("D19").Value = '2.092.44'
$ws.Range("E19").Value = '  +0.19%  '

Set-TextValue "D20" '233.26'
$ws.Range("E20").Value = '  +2.60%  '

Set-TextValue "D21" '12.74'
$ws.Range("E21").Value = '  +0.91%  '

$ws.Range("E22").Value = '  +0.30%  '

Set-TextValue "D23" '7.189'
$ws.Range("E23").Value = '  -0.53%  '

$ws.Range("E24").Value = '  +0.16%  '

Set-TextValue "D25" '160.56'
$ws.Range("E25").Value = '  -0.20%  '

Set-TextValue "D26" '0.1400'
$ws.Range("E26").Value = '  -0.54%  '

Set-TextValue "D27" '8.575'
$ws.Range("E27").Value = '  +0.15%  '

Set-TextValue "D28" '18.00'
$ws.Range("E28").Value = '  -0.06%  '

Set-TextValue "D29" '1.500'
$ws.Range("E29").Value = '  -0.04%  '

Set-TextValue "D30" '4.200'
$ws.Range("E30").Value = '  +1.61%  '

Set-TextValue "D31" '4.164'

Set-TextValue "D32" '0.05580'
$ws.Range("E32").Value = '  +3.00%  '

Set-TextValue "D33" '1.208'
$ws.Range("E33").Value = '  +0.22%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D34" '1.862'
$ws.Range("E34").Value = '  -0.33%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D35" '0.7518'
$ws.Range("E35").Value = '  +0.27%  '

Set-TextValue "D36" '1.151'
$ws.Range("E36").Value = '  +0.63%  '

Set-TextValue "D37" '2.665'
$ws.Range("E37").Value = '  +0.05%  '

$ws.Range("D38").Value = '1.240.65'
$ws.Range("E38").Value = '  -0.09%  '

Set-TextValue "D39" '2.773'
$ws.Range("E39").Value = '  +0.33%  '

Set-TextValue "D40" '0.01796'
$ws.Range("E40").Value = '  +0.23%  '

Set-TextValue "D41" '6.646'
$ws.Range("E41").Value = '  +0.41%  '

Set-TextValue "D42" '0.9039'
$ws.Range("E42").Value = '  -0.30%  '

$ws.Range("E43").Value = '  +0.11%  '

Set-TextValue "D44" '102.54'
$ws.Range("E44").Value = '  +0.20%  '

$ws.Range("D45").Value = '1.993.64'
$ws.Range("E45").Value = '  +0.22%  '

Set-TextValue "D46" '66.76'
$ws.Range("E46").Value = '  +2.34%  '

Set-TextValue "D47" '0.5106'

$ws.Range("E48").Value = '  -3.58%  '

Set-TextValue "D49" '0.4106'
$ws.Range("E49").Value = '  +0.05%  '

Set-TextValue "D50" '9.122'
$ws.Range("E50").Value = '  -0.19%  '
